$d = $word.ActiveDocument

# Remove the existing "_GoBack" bookmark from the 4th paragraph (if present)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# The last paragraph of the body is currently empty; insert the new text there
$paraCount = $d.Paragraphs.Count
$targetPara = $d.Paragraphs.Item($paraCount)
$targetRange = $targetPara.Range
$targetRange.Collapse(0)
$targetRange.Text = "Git创建分支简单又快捷"

# Re-create the "_GoBack" bookmark at the end of that same paragraph
$bookmarkRange = $targetPara.Range
$bookmarkRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
